$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format for rows whose new price would otherwise be
# auto-converted to a number by Excel (values like "1.00", "5.86", "0.0000101").
$textFormatRows = @(5,6,8,9,15,18,19,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textFormatRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '64.220.37'
$ws.Cells.Item(2, 5).Value = '  -0.11%  '
$ws.Cells.Item(3, 4).Value = '3.154.22'
$ws.Cells.Item(3, 5).Value = '  -0.86%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '570.06'
$ws.Cells.Item(5, 5).Value = '  -0.06%  '
$ws.Cells.Item(6, 4).Value = '163.03'
$ws.Cells.Item(6, 5).Value = '  -3.67%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 4).Value = '0.575'
$ws.Cells.Item(8, 5).Value = '  -5.62%  '
$ws.Cells.Item(9, 4).Value = '0.116'
$ws.Cells.Item(9, 5).Value = '  -3.52%  '
$ws.Cells.Item(10, 5).Value = '  -1.70%  '
$ws.Cells.Item(11, 5).Value = '  -1.12%  '
$ws.Cells.Item(12, 4).Value = '3.700.98'
$ws.Cells.Item(12, 5).Value = '  -0.77%  '
$ws.Cells.Item(13, 5).Value = '  -0.84%  '
$ws.Cells.Item(14, 4).Value = '64.296.43'
$ws.Cells.Item(14, 5).Value = '  -0.08%  '
$ws.Cells.Item(15, 4).Value = '25.06'
$ws.Cells.Item(15, 5).Value = '  -1.38%  '
$ws.Cells.Item(16, 4).Value = '3.155.27'
$ws.Cells.Item(16, 5).Value = '  -0.72%  '
$ws.Cells.Item(18, 4).Value = '403.65'
$ws.Cells.Item(18, 5).Value = '  -3.62%  '
$ws.Cells.Item(19, 4).Value = '12.64'
$ws.Cells.Item(19, 5).Value = '  -1.47%  '
$ws.Cells.Item(20, 5).Value = '  -2.52%  '
$ws.Cells.Item(21, 4).Value = '7.09'
$ws.Cells.Item(21, 5).Value = '  +0.13%  '
$ws.Cells.Item(22, 2).Value = 'LEO'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(22, 4).Value = '5.86'
$ws.Cells.Item(22, 5).Value = '  +3.61%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  -0.05%  '
$ws.Cells.Item(24, 2).Value = 'Litecoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).Value = '68.49'
$ws.Cells.Item(24, 5).Value = '  -2.41%  '
$ws.Cells.Item(25, 2).Value = 'Polygon'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(25, 4).Value = '0.482'
$ws.Cells.Item(25, 5).Value = '  -1.81%  '
$ws.Cells.Item(26, 2).Value = 'Kaspa'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(26, 4).Value = '0.193'
$ws.Cells.Item(26, 5).Value = '  -4.65%  '
$ws.Cells.Item(27, 2).Value = 'PEPE'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(27, 4).Value = '0.0000101'
$ws.Cells.Item(27, 5).Value = '  -4.77%  '
$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).Value = '8.80'
$ws.Cells.Item(28, 5).Value = '  +0.02%  '
$ws.Cells.Item(29, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(29, 4).Value = '1.00'
$ws.Cells.Item(29, 5).Value = '  +0.67%  '
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).Value = '1.80'
$ws.Cells.Item(30, 5).Value = '  -1.35%  '
$ws.Cells.Item(31, 2).Value = 'EthereumClassic'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(31, 4).Value = '21.11'
$ws.Cells.Item(31, 5).Value = '  -3.28%  '
$ws.Cells.Item(32, 2).Value = 'Aptos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(32, 4).Value = '6.26'
$ws.Cells.Item(32, 5).Value = '  -1.39%  '
$ws.Cells.Item(33, 2).Value = 'NEARProtocol'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(33, 4).Value = '4.81'
$ws.Cells.Item(33, 5).Value = '  -4.19%  '
$ws.Cells.Item(34, 2).Value = 'Monero'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(34, 4).Value = '156.12'
$ws.Cells.Item(34, 5).Value = '  +0.21%  '
$ws.Cells.Item(35, 2).Value = 'Fetch.AI'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(35, 4).Value = '1.11'
$ws.Cells.Item(35, 5).Value = '  -2.28%  '
$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(36, 4).Value = '1.33'
$ws.Cells.Item(36, 5).Value = '  -3.40%  '
$ws.Cells.Item(37, 2).Value = 'Maker'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(37, 4).Value = '2.666.60'
$ws.Cells.Item(37, 5).Value = '  -1.37%  '
$ws.Cells.Item(38, 2).Value = 'Stacks'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(38, 4).Value = '1.67'
$ws.Cells.Item(38, 5).Value = '  -1.86%  '
$ws.Cells.Item(39, 2).Value = 'EnergySwap'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(39, 4).Value = '23.62'
$ws.Cells.Item(39, 5).Value = '  -4.04%  '
$ws.Cells.Item(40, 2).Value = 'Filecoin'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(40, 4).Value = '4.07'
$ws.Cells.Item(40, 5).Value = '  -2.63%  '
$ws.Cells.Item(41, 2).Value = 'Mantle'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(41, 4).Value = '0.693'
$ws.Cells.Item(41, 5).Value = '  -2.47%  '
$ws.Cells.Item(42, 2).Value = 'Hedera'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(42, 4).Value = '0.0615'
$ws.Cells.Item(42, 5).Value = '  -1.61%  '
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(43, 4).Value = '5.40'
$ws.Cells.Item(43, 5).Value = '  -5.12%  '
$ws.Cells.Item(44, 2).Value = 'VeChain'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(44, 4).Value = '0.0255'
$ws.Cells.Item(44, 5).Value = '  -2.91%  '
$ws.Cells.Item(45, 2).Value = 'Bittensor'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(45, 4).Value = '287.46'
$ws.Cells.Item(45, 5).Value = '  -3.10%  '
$ws.Cells.Item(46, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(46, 4).Value = '21.14'
$ws.Cells.Item(46, 5).Value = '  -3.64%  '
$ws.Cells.Item(47, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(47, 4).Value = '0.999'
$ws.Cells.Item(47, 5).Value = '  -0.01%  '
$ws.Cells.Item(48, 2).Value = 'Stellar'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(48, 4).Value = '0.0977'
$ws.Cells.Item(48, 5).Value = '  -1.70%  '
$ws.Cells.Item(49, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(49, 4).Value = '10.48'
$ws.Cells.Item(49, 5).Value = '  +0.31%  '
$ws.Cells.Item(50, 2).Value = 'dogwifhat'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(50, 4).Value = '1.88'
$ws.Cells.Item(50, 5).Value = '  -8.65%  '
$ws.Cells.Item(51, 2).Value = 'Cosmos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(51, 4).Value = '5.68'
$ws.Cells.Item(51, 5).Value = '  -1.91%  '
